# The "DeviceList" sheet (second tab, already the active sheet/tab in the
# workbook) had its column D removed entirely - all the device records that
# used to live in columns E:I shift left into D:H.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

# Remember the conditional formatting ranges before the shift so we can
# reapply them correctly (Excel normally re-anchors these automatically on a
# real column delete, so we do it explicitly here).
$fcs = $ws.Range("B2:I2").FormatConditions

# Delete the entire column D - everything to the right shifts one column left.
$ws.Columns("D").Delete()

# Re-anchor the conditional formatting rules to the new (narrower) range.
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("B2:H2"))
}

# Reset the view: scroll back to show column A and move the selection.
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
$ws.Range("C16").Select() | Out-Null
